$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-8 (columns A-E)
$data = @(
    @(1, 1, 23, 1380, 45738.05504629629),
    @(2, 2, 2, 140, 45738.05511574074),
    @(3, 1, 15, 900, 45738.14853009259),
    @(4, 1, 15, 900, 45738.14856481482),
    @(5, 1, 2, 120, 45738.14869212963),
    @(6, 1, 15, 900, 45738.15077546296),
    @(7, 1, 20, 1200, 45738.15090277778)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}

# Remove row 9 entirely (was an extra record in the old file)
$ws.Rows.Item(9).Delete()
